# Arabic translation pass for "Email 4-2 [TEMPLATE] Partner email – reminder
# to submit documents.docx". Each call is scoped to the specific paragraph
# (or comment) whose run text changed in the source diff, so that strings
# which are repeated verbatim elsewhere in the document (", ", " or ",
# ", at ", "English", etc.) are only touched where the diff actually shows
# an edit.
#
# NOTE: Range.Find.Execute (when it succeeds) collapses the Range object
# down to the matched/replaced span, so reusing one captured Range across
# several Find/Replace calls on the same paragraph would let later calls
# leak into the rest of the document. To avoid that, every replacement
# re-fetches a brand-new Range for the paragraph/comment right before the
# call.

$d = $word.ActiveDocument

function Replace-InParagraph($index, $find, $replace) {
    $rng = $d.Paragraphs($index).Range
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

function Replace-InComment($index, $find, $replace) {
    $rng = $d.Comments($index).Range
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- Paragraph 1: language picker line (top banner) ---
Replace-InParagraph 1 "English" "الإنجليزية"
Replace-InParagraph 1 " / Portuguese / French / Thai / Vietnamese / Spanish" " /البرتغالية/الفرنسية/التايلندية/الفيتنامية/الإسبانية"

# --- Paragraph 3: "English" language heading ---
Replace-InParagraph 3 "English" "الإنجليزية"

# --- Paragraph 6: Brief description ---
Replace-InParagraph 6 `
    "An email sent to partners in the target country who RSVPed yes but haven’t sent their documents to us. It will be sent via customer.io" `
    "رسالة بريد إلكتروني مرسلة إلى الشركاء في البلد المستهدف الذين قاموا بالرد بنعم ولكنهم لم يرسلوا مستنداتهم إلينا. سيتم إرسالها عبر customer.io"

# --- Paragraph 9: Target audience description ---
Replace-InParagraph 9 `
    "Invited partners who haven’t submitted their documents" `
    "الشركاء المدعوون الذين لم يقدموا مستنداتهم"

# --- Paragraph 14: email #1 heading ---
Replace-InParagraph 14 "Don’t forget to send your documents" "لا تنس إرسال مستنداتك"

# --- Paragraph 16: "Hi [PARTNER NAME], " greeting ---
Replace-InParagraph 16 "Hi " "مرحبًا  "
Replace-InParagraph 16 ", " ",، "

# --- Paragraph 19: "To confirm your registration..." ---
Replace-InParagraph 19 `
    "To confirm your registration, we need the following documents from you by " `
    "لتأكيد تسجيلك، نحتاج إلى المستندات التالية منك بواسطة "

# --- Paragraph 20: bullet placeholder ---
Replace-InParagraph 20 "[insert list of documents required]" "[أدخل قائمة المستندات المطلوبة]"

# --- Paragraph 21: "Please send a copy..." ---
Replace-InParagraph 21 "Please send a copy of these documents to your country manager, " "يرجى إرسال نسخة من هذه المستندات إلى مدير بلدك، "
Replace-InParagraph 21 ", at " "، على "
Replace-InParagraph 21 " or " " أو "
Replace-InParagraph 21 " (WhatsApp), so that we can make the necessary arrangements for you, including accommodation and transportation." " (WhatsApp)، حتى نتمكن من اتخاذ الترتيبات اللازمة لك، بما في ذلك الإقامة والنقل."

# --- Paragraph 22: "If you have any questions, please contact your country manager." ---
Replace-InParagraph 22 `
    "If you have any questions, please contact your country manager." `
    "إذا كانت لديك أي أسئلة، فيُرجى الاتصال بمديرك الإقليمي."

# --- Paragraph 23: sign-off ---
Replace-InParagraph 23 "We look forward to seeing you there!" "نتطلع إلى رؤيتك هناك!"

# --- Paragraph 31: email #2 heading ---
Replace-InParagraph 31 "Don’t forget to send your documents" "لا تنس إرسال مستنداتك"

# --- Paragraph 35: "We’re excited to see you at the upcoming [EVENT NAME]. '" ---
Replace-InParagraph 35 "We’re excited to see you at the upcoming " "نحن متحمسون لرؤيتك في "
Replace-InParagraph 35 ". ‘" " القادم. ‘"

# --- Paragraph 36: "To ensure you have the best experience..." ---
Replace-InParagraph 36 `
    "To ensure you have the best experience at this event, we need the following documents from you by " `
    "لضمان حصولك على أفضل تجربة في هذا الحدث، نحتاج إلى المستندات التالية منك بواسطة "

# --- Paragraph 37: bullet placeholder ---
Replace-InParagraph 37 "[insert list of documents required]" "[أدخل قائمة المستندات المطلوبة]"

# --- Paragraph 38: "Please reply to this email..." ---
Replace-InParagraph 38 `
    "Please reply to this email with a copy of these documents so that we have make the necessary arrangements for you, including accommodation and transportation." `
    "يرجى الرد على هذا البريد الإلكتروني مع نسخة من هذه المستندات حتى نتمكن من اتخاذ الترتيبات اللازمة لك، بما في ذلك الإقامة والنقل."

# --- Paragraph 39: "If you have any questions, please contact us via live chat or WhatsApp." ---
Replace-InParagraph 39 "If you have any questions, please contact us via " "إذا كانت لديك أي أسئلة، فاتصل بنا:  "
Replace-InParagraph 39 "live chat" "الدردشة الحية"

# --- Paragraph 41: sign-off ---
Replace-InParagraph 41 "We look forward to seeing you there!" "نتطلع إلى رؤيتك هناك!"

# --- Comment: "choose either one" ---
Replace-InComment 1 "choose either one" "اختر أيًا منهما"
